$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "'24.856.17"
$ws.Range("E2").Value = "  -3.98%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "'1.679.78"
$ws.Range("E3").Value = "  -3.09%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +2.11%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'309.95"
$ws.Range("E5").Value = "  -0.56%  "

$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "'0.9972"
$ws.Range("E6").Value = "  +1.84%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("D7").Value = "'0.3678"
$ws.Range("E7").Value = "  -2.18%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("D8").Value = "'0.3347"
$ws.Range("E8").Value = "  -5.81%  "

$ws.Range("B9").Value = "OKB"
$ws.Range("D9").Value = "'47.34"
$ws.Range("E9").Value = "  -6.27%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("D10").Value = "'1.176"
$ws.Range("E10").Value = "  -2.88%  "

$ws.Range("B11").Value = "Dogecoin"
$ws.Range("D11").Value = "'0.07343"
$ws.Range("E11").Value = "  -1.90%  "

$ws.Range("B12").Value = "BinanceUSD"
$ws.Range("D12").Value = "'0.9969"
$ws.Range("E12").Value = "  +2.26%  "

$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'20.60"
$ws.Range("E13").Value = "  -4.27%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.181"
$ws.Range("E14").Value = "  -2.04%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("D15").Value = "'6.825"
$ws.Range("E15").Value = "  -1.21%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("D16").Value = "'1.675.71"
$ws.Range("E16").Value = "  -2.43%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("D17").Value = "'0.00001100"
$ws.Range("E17").Value = "  -3.98%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("D18").Value = "'0.06605"
$ws.Range("E18").Value = "  -1.61%  "

$ws.Range("B19").Value = "Dai"
$ws.Range("D19").Value = "'0.9972"
$ws.Range("E19").Value = "  +2.05%  "

$ws.Range("B20").Value = "Litecoin"
$ws.Range("D20").Value = "'82.21"
$ws.Range("E20").Value = "  -4.19%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("D21").Value = "'16.88"
$ws.Range("E21").Value = "  -1.10%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("D22").Value = "'6.189"
$ws.Range("E22").Value = "  -1.45%  "

$ws.Range("B23").Value = "Cosmos"
$ws.Range("D23").Value = "'12.58"
$ws.Range("E23").Value = "  +1.36%  "

$ws.Range("B24").Value = "WrappedBTC"
$ws.Range("D24").Value = "'24.787.43"
$ws.Range("E24").Value = "  -3.39%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("D25").Value = "'2.429"
$ws.Range("E25").Value = "  +1.28%  "

$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("D26").Value = "'2.723"
$ws.Range("E26").Value = "  -4.47%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("D27").Value = "'19.84"
$ws.Range("E27").Value = "  -2.03%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("D28").Value = "'149.78"
$ws.Range("E28").Value = "  -3.12%  "

$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("D29").Value = "'130.30"
$ws.Range("E29").Value = "  -0.47%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.235"
$ws.Range("E30").Value = "  +7.98%  "

$ws.Range("B31").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C31").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D31").Value = "'1.858.38"
$ws.Range("E31").Value = "  -2.27%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("D32").Value = "'6.536"
$ws.Range("E32").Value = "  -1.20%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("D33").Value = "'4.144"
$ws.Range("E33").Value = "  +1.62%  "

$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'13.49"
$ws.Range("E34").Value = "  +1.82%  "

$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'1.739"
$ws.Range("E35").Value = "  -2.59%  "

$ws.Range("B36").Value = "Stellar"
$ws.Range("D36").Value = "'0.08607"
$ws.Range("E36").Value = "  +1.16%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("D37").Value = "'5.459"
$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("D38").Value = "'0.06482"
$ws.Range("E38").Value = "  -1.88%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02349"
$ws.Range("E39").Value = "  -2.33%  "

$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.777"
$ws.Range("E40").Value = "  -4.26%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("D41").Value = "'0.2172"
$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("D42").Value = "'1.247"
$ws.Range("E42").Value = "  -1.22%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("D43").Value = "'0.6285"
$ws.Range("E43").Value = "  -1.51%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("D44").Value = "'0.9962"
$ws.Range("E44").Value = "  +1.85%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("D45").Value = "'13.46"
$ws.Range("E45").Value = "  -0.81%  "

$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("D46").Value = "'3.798"
$ws.Range("E46").Value = "  -1.27%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("D47").Value = "'0.5982"
$ws.Range("E47").Value = "  -3.30%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("D48").Value = "'2.048"
$ws.Range("E48").Value = "  -3.14%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("D49").Value = "'125.60"
$ws.Range("E49").Value = "  -4.26%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("D50").Value = "'0.07158"
$ws.Range("E50").Value = "  -4.25%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("D51").Value = "'77.32"
$ws.Range("E51").Value = "  -0.82%  "
